$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comparison")

# Update projection results with new auto-enrollment and auto-increase results
# Row 2
$ws.Cells.Item(2, 3).Value = 9478
$ws.Cells.Item(2, 4).Value = 9456
$ws.Cells.Item(2, 5).Value = 8369
$ws.Cells.Item(2, 6).Value = 0.8850465313028765
$ws.Cells.Item(2, 7).Value = 0.8829921924456636
$ws.Cells.Item(2, 8).Value = 0.09597596898787557
$ws.Cells.Item(2, 9).Value = 0.08474603127870128
$ws.Cells.Item(2, 10).Value = 41011357.62888187
$ws.Cells.Item(2, 11).Value = 14306014.57377693
$ws.Cells.Item(2, 13).Value = 14306014.57377693
$ws.Cells.Item(2, 14).Value = 55317372.20265881
$ws.Cells.Item(2, 15).Value = 812856555.5672001
$ws.Cells.Item(2, 16).Value = 793677941.7132001
$ws.Cells.Item(2, 17).Value = 0.01759967915100887
$ws.Cells.Item(2, 18).Value = 0.01802496179104659

# Row 3
$ws.Cells.Item(3, 3).Value = 9762
$ws.Cells.Item(3, 4).Value = 9743
$ws.Cells.Item(3, 5).Value = 8661
$ws.Cells.Item(3, 6).Value = 0.8889459098840193
$ws.Cells.Item(3, 7).Value = 0.8872157344806392
$ws.Cells.Item(3, 8).Value = 0.09352092230585833
$ws.Cells.Item(3, 9).Value = 0.08297323377289889
$ws.Cells.Item(3, 10).Value = 43255456.65684056
$ws.Cells.Item(3, 11).Value = 15158002.11521988
$ws.Cells.Item(3, 13).Value = 15158002.11521988
$ws.Cells.Item(3, 14).Value = 58413458.77206045
$ws.Cells.Item(3, 15).Value = 868756639.5402131
$ws.Cells.Item(3, 16).Value = 848873415.5163431
$ws.Cells.Item(3, 17).Value = 0.01744792664058626
$ws.Cells.Item(3, 18).Value = 0.01785661070090143

# Row 4
$ws.Cells.Item(4, 3).Value = 10048
$ws.Cells.Item(4, 4).Value = 10026
$ws.Cells.Item(4, 5).Value = 8905
$ws.Cells.Item(4, 6).Value = 0.8881907041691602
$ws.Cells.Item(4, 7).Value = 0.8862460191082803
$ws.Cells.Item(4, 8).Value = 0.09177983476932886
$ws.Cells.Item(4, 9).Value = 0.08133951319873342
$ws.Cells.Item(4, 10).Value = 45461625.99696768
$ws.Cells.Item(4, 11).Value = 15941015.34003143
$ws.Cells.Item(4, 13).Value = 15941015.34003143
$ws.Cells.Item(4, 14).Value = 61402641.33699911
$ws.Cells.Item(4, 15).Value = 923597016.4015658
$ws.Cells.Item(4, 16).Value = 902839601.727118
$ws.Cells.Item(4, 17).Value = 0.01725970857088664
$ws.Cells.Item(4, 18).Value = 0.01765653091594179

# Row 5
$ws.Cells.Item(5, 3).Value = 10340
$ws.Cells.Item(5, 4).Value = 10303
$ws.Cells.Item(5, 5).Value = 9132
$ws.Cells.Item(5, 6).Value = 0.8863437833640687
$ws.Cells.Item(5, 7).Value = 0.8831721470019342
$ws.Cells.Item(5, 8).Value = 0.09026420469662112
$ws.Cells.Item(5, 9).Value = 0.07971883145933695
$ws.Cells.Item(5, 10).Value = 47924844.26765846
$ws.Cells.Item(5, 11).Value = 16839396.97316202
$ws.Cells.Item(5, 13).Value = 16839396.97316202
$ws.Cells.Item(5, 14).Value = 64764241.24082048
$ws.Cells.Item(5, 15).Value = 981614581.1341684
$ws.Cells.Item(5, 16).Value = 958614956.0512464
$ws.Cells.Item(5, 17).Value = 0.01715479506600808
$ws.Cells.Item(5, 18).Value = 0.01756638248429519

# Row 6
$ws.Cells.Item(6, 3).Value = 10667
$ws.Cells.Item(6, 4).Value = 10640
$ws.Cells.Item(6, 5).Value = 9432
$ws.Cells.Item(6, 6).Value = 0.8864661654135338
$ws.Cells.Item(6, 7).Value = 0.8842223680509984
$ws.Cells.Item(6, 8).Value = 0.08827800578727582
$ws.Cells.Item(6, 9).Value = 0.07805738732404476
$ws.Cells.Item(6, 10).Value = 50668659.16144493
$ws.Cells.Item(6, 11).Value = 17839816.11902389
$ws.Cells.Item(6, 13).Value = 17839816.11902389
$ws.Cells.Item(6, 14).Value = 68508475.28046882
$ws.Cells.Item(6, 15).Value = 1049528624.192984
$ws.Cells.Item(6, 16).Value = 1025805590.179856
$ws.Cells.Item(6, 17).Value = 0.01699793193610273
$ws.Cells.Item(6, 18).Value = 0.01739103031783636

# Row 7
$ws.Cells.Item(7, 3).Value = 9478
$ws.Cells.Item(7, 4).Value = 9454
$ws.Cells.Item(7, 5).Value = 8552
$ws.Cells.Item(7, 6).Value = 0.9045906494605458
$ws.Cells.Item(7, 7).Value = 0.9023000633044946
$ws.Cells.Item(7, 8).Value = 0.09509490102641401
$ws.Cells.Item(7, 9).Value = 0.08582945701391566
$ws.Cells.Item(7, 10).Value = 41871748.61682985
$ws.Cells.Item(7, 11).Value = 14736210.06775092
$ws.Cells.Item(7, 13).Value = 14736210.06775092
$ws.Cells.Item(7, 14).Value = 56607958.68458077
$ws.Cells.Item(7, 15).Value = 817193256.2772001
$ws.Cells.Item(7, 16).Value = 797814518.5032001
$ws.Cells.Item(7, 17).Value = 0.01803271130122035
$ws.Cells.Item(7, 18).Value = 0.0184707218607627

# Row 8
$ws.Cells.Item(8, 3).Value = 9764
$ws.Cells.Item(8, 4).Value = 9740
$ws.Cells.Item(8, 5).Value = 8841
$ws.Cells.Item(8, 6).Value = 0.907700205338809
$ws.Cells.Item(8, 7).Value = 0.9054690700532568
$ws.Cells.Item(8, 8).Value = 0.09318518303774719
$ws.Cells.Item(8, 9).Value = 0.08440088111805848
$ws.Cells.Item(8, 10).Value = 44278388.44033591
$ws.Cells.Item(8, 11).Value = 15669468.00696756
$ws.Cells.Item(8, 13).Value = 15669468.00696756
$ws.Cells.Item(8, 14).Value = 59947856.44730347
$ws.Cells.Item(8, 15).Value = 866217456.3869131
$ws.Cells.Item(8, 16).Value = 846193941.7348431
$ws.Cells.Item(8, 17).Value = 0.01808953155057231
$ws.Cells.Item(8, 18).Value = 0.01851758472158575

# Row 9
$ws.Cells.Item(9, 3).Value = 10042
$ws.Cells.Item(9, 4).Value = 10027
$ws.Cells.Item(9, 5).Value = 9131
$ws.Cells.Item(9, 6).Value = 0.910641268574848
$ws.Cells.Item(9, 7).Value = 0.9092810197171878
$ws.Cells.Item(9, 8).Value = 0.0918294128663809
$ws.Cells.Item(9, 9).Value = 0.08351367943466678
$ws.Cells.Item(9, 10).Value = 47517624.66176366
$ws.Cells.Item(9, 11).Value = 16969014.67242942
$ws.Cells.Item(9, 13).Value = 16969014.67242942
$ws.Cells.Item(9, 14).Value = 64486639.33419308
$ws.Cells.Item(9, 15).Value = 926735018.4242668
$ws.Cells.Item(9, 16).Value = 906466055.975973
$ws.Cells.Item(9, 17).Value = 0.01831053573575102
$ws.Cells.Item(9, 18).Value = 0.0187199670197901

# Row 10
$ws.Cells.Item(10, 3).Value = 10325
$ws.Cells.Item(10, 4).Value = 10298
$ws.Cells.Item(10, 5).Value = 9404
$ws.Cells.Item(10, 6).Value = 0.9131870266071082
$ws.Cells.Item(10, 7).Value = 0.9107990314769976
$ws.Cells.Item(10, 8).Value = 0.09067876875608344
$ws.Cells.Item(10, 9).Value = 0.08261628487963281
$ws.Cells.Item(10, 10).Value = 50835217.98122857
$ws.Cells.Item(10, 11).Value = 18294583.82994707
$ws.Cells.Item(10, 13).Value = 18294583.82994707
$ws.Cells.Item(10, 14).Value = 69129801.81117564
$ws.Cells.Item(10, 15).Value = 982952967.4198503
$ws.Cells.Item(10, 16).Value = 962400924.913815
$ws.Cells.Item(10, 17).Value = 0.01861186082785676
$ws.Cells.Item(10, 18).Value = 0.01900931655025725

# Row 11
$ws.Cells.Item(11, 3).Value = 10639
$ws.Cells.Item(11, 4).Value = 10615
$ws.Cells.Item(11, 5).Value = 9724
$ws.Cells.Item(11, 6).Value = 0.9160621761658031
$ws.Cells.Item(11, 7).Value = 0.9139956762853652
$ws.Cells.Item(11, 8).Value = 0.08959285875955238
$ws.Cells.Item(11, 9).Value = 0.08191004404341455
$ws.Cells.Item(11, 10).Value = 54617108.1983126
$ws.Cells.Item(11, 11).Value = 19814040.63745773
$ws.Cells.Item(11, 13).Value = 19814040.63745773
$ws.Cells.Item(11, 14).Value = 74431148.83577034
$ws.Cells.Item(11, 15).Value = 1044291442.661237
$ws.Cells.Item(11, 16).Value = 1022775193.783701
$ws.Cells.Item(11, 17).Value = 0.01897366944515441
$ws.Cells.Item(11, 18).Value = 0.01937282088760557

# Row 12
$ws.Cells.Item(12, 3).Value = 9478
$ws.Cells.Item(12, 4).Value = 9456
$ws.Cells.Item(12, 5).Value = 9456
$ws.Cells.Item(12, 7).Value = 0.9976788351972991
$ws.Cells.Item(12, 8).Value = 0.09258994317483452
$ws.Cells.Item(12, 9).Value = 0.092375026657653
$ws.Cells.Item(12, 10).Value = 47112440.0389384
$ws.Cells.Item(12, 11).Value = 17356555.7788052
$ws.Cells.Item(12, 13).Value = 17356555.7788052
$ws.Cells.Item(12, 14).Value = 64468995.81774361
$ws.Cells.Item(12, 15).Value = 813109489.4172001
$ws.Cells.Item(12, 16).Value = 794833158.7732
$ws.Cells.Item(12, 17).Value = 0.02134590237194943
$ws.Cells.Item(12, 18).Value = 0.02183672835893573

# Row 13
$ws.Cells.Item(13, 3).Value = 9762
$ws.Cells.Item(13, 4).Value = 9738
$ws.Cells.Item(13, 5).Value = 9738
$ws.Cells.Item(13, 7).Value = 0.9975414874001229
$ws.Cells.Item(13, 8).Value = 0.09761071208907925
$ws.Cells.Item(13, 9).Value = 0.09737073492352527
$ws.Cells.Item(13, 10).Value = 55505811.34755692
$ws.Cells.Item(13, 11).Value = 21282186.73698258
$ws.Cells.Item(13, 13).Value = 21282186.73698258
$ws.Cells.Item(13, 14).Value = 76787998.0845395
$ws.Cells.Item(13, 15).Value = 863794315.9995871
$ws.Cells.Item(13, 16).Value = 845379885.828617
$ws.Cells.Item(13, 17).Value = 0.0246380259082334
$ws.Cells.Item(13, 18).Value = 0.0251747020407546

# Row 14
$ws.Cells.Item(14, 3).Value = 10046
$ws.Cells.Item(14, 4).Value = 10024
$ws.Cells.Item(14, 5).Value = 10024
$ws.Cells.Item(14, 7).Value = 0.9978100736611587
$ws.Cells.Item(14, 8).Value = 0.1008389917673471
$ws.Cells.Item(14, 9).Value = 0.1006181618032936
$ws.Cells.Item(14, 10).Value = 62994171.20052955
$ws.Cells.Item(14, 11).Value = 24707287.94181236
$ws.Cells.Item(14, 13).Value = 24707287.94181236
$ws.Cells.Item(14, 14).Value = 87701459.1423419
$ws.Cells.Item(14, 15).Value = 916693029.3738154
$ws.Cells.Item(14, 16).Value = 897970254.2003546
$ws.Cells.Item(14, 17).Value = 0.02695262988820771
$ws.Cells.Item(14, 18).Value = 0.02751459508401454

# Row 15
$ws.Cells.Item(15, 3).Value = 10337
$ws.Cells.Item(15, 4).Value = 10309
$ws.Cells.Item(15, 5).Value = 10309
$ws.Cells.Item(15, 7).Value = 0.9972912837380284
$ws.Cells.Item(15, 8).Value = 0.1007746111987947
$ws.Cells.Item(15, 9).Value = 0.1005016413706467
$ws.Cells.Item(15, 10).Value = 67918402.13472392
$ws.Cells.Item(15, 11).Value = 26836175.90669475
$ws.Cells.Item(15, 13).Value = 26836175.90669475
$ws.Cells.Item(15, 14).Value = 94754578.04141869
$ws.Cells.Item(15, 15).Value = 969816544.9170408
$ws.Cells.Item(15, 16).Value = 951407665.5136355
$ws.Cells.Item(15, 17).Value = 0.02767139419032117
$ws.Cells.Item(15, 18).Value = 0.02820681068635992

# Row 16
$ws.Cells.Item(16, 3).Value = 10656
$ws.Cells.Item(16, 4).Value = 10631
$ws.Cells.Item(16, 5).Value = 10631
$ws.Cells.Item(16, 7).Value = 0.9976539039039038
$ws.Cells.Item(16, 8).Value = 0.1002504979019584
$ws.Cells.Item(16, 9).Value = 0.1000153006001989
$ws.Cells.Item(16, 10).Value = 73124236.81502223
$ws.Cells.Item(16, 11).Value = 29067604.94581254
$ws.Cells.Item(16, 13).Value = 29067604.94581254
$ws.Cells.Item(16, 14).Value = 102191841.7608348
$ws.Cells.Item(16, 15).Value = 1033114253.054143
$ws.Cells.Item(16, 16).Value = 1014594770.022316
$ws.Cells.Item(16, 17).Value = 0.02813590545274296
$ws.Cells.Item(16, 18).Value = 0.02864947248365294
